$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 161, pushing existing rows 161-234 down to 162-235
$ws.Rows.Item(161).Insert()

# Populate the new row 161 with the new data record
$ws.Range("A161").Value = 11
$ws.Range("B161").Value = "Vega Monumental Concepción"
$ws.Range("C161").Value = "Bíobío"
$ws.Range("D161").Value = 44845
$ws.Range("E161").Value = 8
$ws.Range("F161").Value = 100112040
$ws.Range("G161").Value = "Cilantro"
$ws.Range("H161").Value = "Sin especificar"
$ws.Range("I161").Value = "Primera"
$ws.Range("J161").Value = 150
$ws.Range("K161").Value = 6000
$ws.Range("L161").Value = 6500
$ws.Range("M161").Value = 6333
$ws.Range("N161").Value = "$/caja 36 atados"
$ws.Range("O161").Value = "Región Metropolitana"
$ws.Range("P161").Value = 176
$ws.Range("Q161").Value = 36
$ws.Range("R161").Value = "Hortaliza"
